$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Update the full product-name cell (B1) on both sheets to the new, shortened text.
$ws1.Range("B1").Value = "4256-MS-EI-DB-DL-REC-CTRFD-RNI-FEE+INT-FFC-SAR-FFROP-DAILY-FIFR-1-MD-TR-1-ONT-PE-1st"
$ws2.Range("B1").Value = "4256-MS-EI-DB-DL-REC-CTRFD-RNI-FEE+INT-FFC-SAR-FFROP-DAILY-FIFR-1-MD-TR-1-ONT-PE-1st"

# Update the short name (B2) on ProductLoanInput to the new text value.
$ws1.Range("B2").Value = "425z"

# Change the active sheet / selection: ProductLoanInput loses the tab selection and its
# selection moves from B17 to B3; ProductLoanOutput becomes the selected tab.
$ws1.Activate()
$ws1.Range("B3").Select()
$ws2.Activate()
$ws2.Range("B1").Select()
